$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2007, 2008 and 2009 data rows (rows 2-4) are removed entirely; the
# remaining years (2010-2013, previously rows 5-8) shift up to become
# rows 2-5, and the sheet's used range shrinks from A1:F8 to A1:F5.
$ws.Range("A2:A4").EntireRow.Delete()
